$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "Notes" + <bookmark _GoBack/> + " taker: " -> single run
# "Notes taker: " (bookmark removed).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Notes taker: ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Notes taker: ", 2) | Out-Null

# ------------------------------------------------------------------
# Change 2: append " drafts" after "... C4 architecture model" and
# move the _GoBack bookmark to sit right after the new text.
# ------------------------------------------------------------------

# Drop the (now stale) _GoBack bookmark wherever it currently lives;
# we'll recreate it at the correct spot once the new text is in.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$target = $d.Content.Find
$target.ClearFormatting()
$found = $d.Content.Find.Execute("C4 architecture model", $true)

$hit = $d.Content
$hit.Find.ClearFormatting()
$hit.Find.Execute("C4 architecture model", $true) | Out-Null
$insertAt = $hit.End

$ins = $d.Range($insertAt, $insertAt)
$ins.InsertAfter(" drafts")

# Re-select just the freshly inserted " drafts" text and nudge a
# formatting property on/off so it is kept in its own run instead of
# being re-merged with the preceding run.
$newText = $d.Range($insertAt, $insertAt + 7)
$newText.Font.Bold = 1
$newText.Font.Bold = 0

# Re-add the _GoBack bookmark collapsed right at the end of the
# paragraph, after " drafts".
$para = $newText.Paragraphs(1).Range
$endPos = $para.End - 1

# Placing a collapsed bookmark exactly at end-of-paragraph-minus-one
# lands it at the wrong spot, so park a throwaway character there
# first, add the bookmark right before it, then remove the throwaway
# character again.
$guard = $d.Range($endPos, $endPos)
$guard.InsertAfter("Z")
$bmSpot = $d.Range($endPos, $endPos)
$bm = $d.Bookmarks.Add("_GoBack", $bmSpot)
$guardRange = $d.Range($endPos, $endPos + 1)
$guardRange.Delete()
